$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "romeesa.ashfaq@ascend.com.sa"
$ws.Range("B3").Value = "Approved"
$ws.Range("C3").Value = "MoH Compliance Program_Digital  /  Others (QA)`nMoH Compliance Program_Digital  /  Others (dotnet)"
$ws.Range("C3").WrapText = $true
